{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Change summary (from the diff):\n//   1. In the \"File to test the 3 designs work\" section, the fragment\n//        \"type of point it wants to create\"\n//      becomes\n//        \"design they want to use and to enter the coordinates\"\n//      (i.e. the sentence \"Running this file will prompt the user for the\n//      type of point it wants to create. It will then run tests ...\"\n//      becomes \"Running this file will prompt the user for the design they\n//      want to use and to enter the coordinates. It will then run tests\n//      ...\"), with a \"_GoBack\" bookmark landing inside the freshly typed\n//      text, between \"coordina\" and \"tes\" (i.e. inside \"coordinates\") -\n//      exactly where the author's cursor was when they stopped editing.\n//   2. The \"_GoBack\" bookmark that used to sit at the end of the\n//      \"Comparison to hypotheses from E.26\" paragraph is gone - Word only\n//      ever keeps a single \"_GoBack\" (the last edit position), so it\n//      simply moved from there to the spot above.\n\n// 1. Drop the old \"_GoBack\" first so only one ever exists at a time.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Find the exact phrase being replaced (scoped tightly so only the\n//    words that actually change are touched; the rest of the sentence -\n//    \"Running this file will prompt the user for the \" / \". It will then\n//    run tests on methods to ensure that the designs work. ...\" - is left\n//    alone).\nconst body = context.document.body;\nconst hits = body.search(\"type of point it wants to create\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nconst target = hits.items[0];\ntarget.insertText(\"design they want to use and to enter the coordinates\", \"Replace\");\nawait context.sync();\n\n// 3. Re-plant \"_GoBack\" inside that new text, right after \"coordina\" (i.e.\n//    between \"coordina\" and \"tes\" in \"coordinates\"). Scope the search to\n//    the paragraph we just edited so the pre-existing, unrelated\n//    \"coordinates\" elsewhere in the document isn't matched instead.\nconst paragraph = target.paragraphs.getFirst();\nconst coordHits = paragraph.search(\"coordina\", { matchCase: true });\ncoordHits.load(\"items\");\nawait context.sync();\n\nconst afterCoordina = coordHits.items[0].getRange(\"End\");\nafterCoordina.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Change summary (from the diff):\n#   1. In the \"File to test the 3 designs work\" section, the fragment\n#        \"type of point it wants to create\"\n#      becomes\n#        \"design they want to use and to enter the coordinates\"\n#      (i.e. the sentence \"Running this file will prompt the user for the\n#      type of point it wants to create. It will then run tests ...\"\n#      becomes \"Running this file will prompt the user for the design they\n#      want to use and to enter the coordinates. It will then run tests\n#      ...\"), with a \"_GoBack\" bookmark landing inside the freshly typed\n#      text, between \"coordina\" and \"tes\" (i.e. inside \"coordinates\") -\n#      exactly where the author's cursor was when they stopped editing.\n#   2. The \"_GoBack\" bookmark that used to sit at the end of the\n#      \"Comparison to hypotheses from E.26\" paragraph is gone - Word only\n#      ever keeps a single \"_GoBack\" (the last edit position), so it\n#      simply moved from there to the spot above.\n\n$d = $word.ActiveDocument\n\n# 1. Drop the old \"_GoBack\" first so only one ever exists at a time.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Find-and-replace the exact phrase being changed (leaves the rest of\n#    the sentence - \"Running this file will prompt the user for the \" /\n#    \". It will then run tests on methods to ensure that the designs\n#    work. ...\" - untouched).\n$findRange = $d.Content\n$oldText = \"type of point it wants to create\"\n$newText = \"design they want to use and to enter the coordinates\"\n$replaced = $findRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n# 3. Re-plant \"_GoBack\" inside that new text, right after \"coordina\" (i.e.\n#    between \"coordina\" and \"tes\" in \"coordinates\"). Scope the second\n#    search to the paragraph we just edited (via its own Range) so the\n#    pre-existing, unrelated \"coordinates\" elsewhere in the document can't\n#    be matched instead.\n$paragraphRange = $findRange.Paragraphs.Item(1).Range\n$paragraphRange.Find.Execute(\"coordina\") | Out-Null\n\n$bookmarkSpot = $d.Range($paragraphRange.End, $paragraphRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot)\n"}
